$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header: "Visit" -> "visit" (lowercase) on sheet1 only
$ws.Range("C1").Value = "visit"

# Column C rows 2-32: numeric 3 -> text "V3"
for ($r = 2; $r -le 32; $r++) {
    $ws.Cells.Item($r, 3).Value = "V3"
}

# Column C rows 33-63: numeric 4 -> text "V4"
for ($r = 33; $r -le 63; $r++) {
    $ws.Cells.Item($r, 3).Value = "V4"
}
